$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$rng = $p1.Range
$frag = @'
<w:p><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251663360" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="0FD85AE0" wp14:editId="67F4FC18"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>220600</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>3359176</wp:posOffset></wp:positionV><wp:extent cx="491490" cy="261824"/><wp:effectExtent l="57150" t="38100" r="60960" b="81280"/><wp:wrapNone/><wp:docPr id="4" name="Cuadro de texto 4"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="491490" cy="261824"/></a:xfrm><a:prstGeom prst="roundRect"><a:avLst/></a:prstGeom><a:solidFill><a:schemeClr val="accent4"/></a:solidFill><a:ln/></wps:spPr><wps:style><a:lnRef idx="0"><a:schemeClr val="accent4"/></a:lnRef><a:fillRef idx="3"><a:schemeClr val="accent4"/></a:fillRef><a:effectRef idx="3"><a:schemeClr val="accent4"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:txbx><w:txbxContent><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="44"/><w:szCs w:val="40"/><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:roundrect w14:anchorId="0FD85AE0" id="Cuadro de texto 4" o:spid="_x0000_s1026" style="position:absolute;margin-left:17.35pt;margin-top:264.5pt;width:38.7pt;height:20.6pt;z-index:251663360;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:top" arcsize="10923f" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQCUNVKQgwIAAHIFAAAOAAAAZHJzL2Uyb0RvYy54bWysVElvGyEUvlfqf0Dcm/GWzco4ch2lqhQl&#xA;UZwqZ8xAjAQ8Ctgz7q/vg1mSpumhVS8zwPve9r3l4rIxmuyFDwpsScdHI0qE5VAp+1zSb4/Xn84o&#xA;CZHZimmwoqQHEejl4uOHi9rNxQS2oCvhCRqxYV67km5jdPOiCHwrDAtH4IRFoQRvWMSrfy4qz2q0&#xA;bnQxGY1Oihp85TxwEQK+XrVCusj2pRQ83kkZRCS6pBhbzF+fv5v0LRYXbP7smdsq3oXB/iEKw5RF&#xA;p4OpKxYZ2Xn1mymjuIcAMh5xMAVIqbjIOWA249GbbNZb5kTOBckJbqAp/D+z/HZ/74mqSjqlxDKD&#xA;JVrtWOWBVIJE0UQg00RS7cIcsWuH6Nh8hgaL3b8HfEy5N9Kb9MesCMqR7sNAMVoiPClNppPjE/TF&#xA;UTY7O51Oz5OZ4kXb+RC/CDAkHUrqYWerB6xjppftb0Js8T0ueQygVXWttM6X1DtipT3ZM6w641zY&#xA;mINFL78gtU2mUmptCvkUD1okM9o+CInM5Ez+aHfWRZ/RCSUxikFxmoPOzfw2oF6xwydVkfv1b5QH&#xA;jewZbByUjbLg3/OuBypki+8ZaPNOFMRm03Ql30B1wIp7aAcnOH6tsCw3LMR75nFSsMg4/fEOP1JD&#xA;XVLoTpRswf947z3hsYFRSkmNk1fS8H3HvKBEf7XY2ufj2SyNar7Mjk8nePGvJZvXErszK8Ayj3HP&#xA;OJ6PCR91f5QezBMuiWXyiiJmOfouaeyPq9juA1wyXCyXGYTD6Vi8sWvHk+lEb+q3x+aJedd1ZpqO&#xA;W+hnlM3f9GaLTZoWlrsIUuXGTQS3rHbE42Dn/u+WUNocr+8Z9bIqFz8BAAD//wMAUEsDBBQABgAI&#xA;AAAAIQBjz3Xb4QAAAAsBAAAPAAAAZHJzL2Rvd25yZXYueG1sTI/BTsMwDIbvSLxDZCRuLGnZuqo0&#xA;nRjSLiAkGLvs5jWhjWiSrsna8vaYE9xs+dfn7y83s+3YqIdgvJOQLAQw7WqvjGskHD52dzmwENEp&#xA;7LzTEr51gE11fVViofzk3vW4jw0jiAsFSmhj7AvOQ91qi2Hhe+3o9ukHi5HWoeFqwIngtuOpEBm3&#xA;aBx9aLHXT62uv/YXS5SDyczxvHt9G9NnMW1FfMFjlPL2Zn58ABb1HP/C8KtP6lCR08lfnAqsk3Cf&#xA;p2uK0pCtEmCUWK0zKnOSsFwmOfCq5P87VD8AAAD//wMAUEsBAi0AFAAGAAgAAAAhALaDOJL+AAAA&#xA;4QEAABMAAAAAAAAAAAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEA&#xA;OP0h/9YAAACUAQAACwAAAAAAAAAAAAAAAAAvAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEA&#xA;lDVSkIMCAAByBQAADgAAAAAAAAAAAAAAAAAuAgAAZHJzL2Uyb0RvYy54bWxQSwECLQAUAAYACAAA&#xA;ACEAY8912+EAAAALAQAADwAAAAAAAAAAAAAAAADdBAAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAAE&#xA;AAQA8wAAAOsFAAAAAA==&#xA;" fillcolor="#ffc000 [3207]" stroked="f"><v:shadow on="t" color="black" opacity="41287f" offset="0,1.5pt"/><v:textbox><w:txbxContent><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="44"/><w:szCs w:val="40"/><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p></w:txbxContent></v:textbox></v:roundrect></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251661312" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="331D99BF" wp14:editId="4CEA94F2"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>2429899</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>2318233</wp:posOffset></wp:positionV><wp:extent cx="1232563" cy="487339"/><wp:effectExtent l="57150" t="38100" r="62865" b="84455"/><wp:wrapNone/><wp:docPr id="3" name="Cuadro de texto 3"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1232563" cy="487339"/></a:xfrm><a:prstGeom prst="roundRect"><a:avLst/></a:prstGeom><a:solidFill><a:schemeClr val="accent1"/></a:solidFill><a:ln/></wps:spPr><wps:style><a:lnRef idx="0"><a:schemeClr val="accent4"/></a:lnRef><a:fillRef idx="3"><a:schemeClr val="accent4"/></a:fillRef><a:effectRef idx="3"><a:schemeClr val="accent4"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:txbx><w:txbxContent><w:p w14:paraId="70180F89" w14:textId="3DF3CBB2" w:rsidR="00347169" w:rsidRPr="00347169" w:rsidRDefault="00347169" w:rsidP="00347169"><w:pPr><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="44"/><w:szCs w:val="40"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r w:rsidRPr="00347169"><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="44"/><w:szCs w:val="40"/><w:lang w:val="es-ES"/></w:rPr><w:t>AYUDA</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:roundrect w14:anchorId="331D99BF" id="Cuadro de texto 3" o:spid="_x0000_s1027" style="position:absolute;margin-left:191.35pt;margin-top:182.55pt;width:97.05pt;height:38.35pt;z-index:251661312;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:top" arcsize="10923f" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQCUNVKQgwIAAHIFAAAOAAAAZHJzL2Uyb0RvYy54bWysVElvGyEUvlfqf0Dcm/GWzco4ch2lqhQl&#xA;UZwqZ8xAjAQ8Ctgz7q/vg1mSpumhVS8zwPve9r3l4rIxmuyFDwpsScdHI0qE5VAp+1zSb4/Xn84o&#xA;CZHZimmwoqQHEejl4uOHi9rNxQS2oCvhCRqxYV67km5jdPOiCHwrDAtH4IRFoQRvWMSrfy4qz2q0&#xA;bnQxGY1Oihp85TxwEQK+XrVCusj2pRQ83kkZRCS6pBhbzF+fv5v0LRYXbP7smdsq3oXB/iEKw5RF&#xA;p4OpKxYZ2Xn1mymjuIcAMh5xMAVIqbjIOWA249GbbNZb5kTOBckJbqAp/D+z/HZ/74mqSjqlxDKD&#xA;JVrtWOWBVIJE0UQg00RS7cIcsWuH6Nh8hgaL3b8HfEy5N9Kb9MesCMqR7sNAMVoiPClNppPjE/TF&#xA;UTY7O51Oz5OZ4kXb+RC/CDAkHUrqYWerB6xjppftb0Js8T0ueQygVXWttM6X1DtipT3ZM6w641zY&#xA;mINFL78gtU2mUmptCvkUD1okM9o+CInM5Ez+aHfWRZ/RCSUxikFxmoPOzfw2oF6xwydVkfv1b5QH&#xA;jewZbByUjbLg3/OuBypki+8ZaPNOFMRm03Ql30B1wIp7aAcnOH6tsCw3LMR75nFSsMg4/fEOP1JD&#xA;XVLoTpRswf947z3hsYFRSkmNk1fS8H3HvKBEf7XY2ufj2SyNar7Mjk8nePGvJZvXErszK8Ayj3HP&#xA;OJ6PCR91f5QezBMuiWXyiiJmOfouaeyPq9juA1wyXCyXGYTD6Vi8sWvHk+lEb+q3x+aJedd1ZpqO&#xA;W+hnlM3f9GaLTZoWlrsIUuXGTQS3rHbE42Dn/u+WUNocr+8Z9bIqFz8BAAD//wMAUEsDBBQABgAI&#xA;AAAAIQBjz3Xb4QAAAAsBAAAPAAAAZHJzL2Rvd25yZXYueG1sTI/BTsMwDIbvSLxDZCRuLGnZuqo0&#xA;nRjSLiAkGLvs5jWhjWiSrsna8vaYE9xs+dfn7y83s+3YqIdgvJOQLAQw7WqvjGskHD52dzmwENEp&#xA;7LzTEr51gE11fVViofzk3vW4jw0jiAsFSmhj7AvOQ91qi2Hhe+3o9ukHi5HWoeFqwIngtuOpEBm3&#xA;aBx9aLHXT62uv/YXS5SDyczxvHt9G9NnMW1FfMFjlPL2Zn58ABb1HP/C8KtP6lCR08lfnAqsk3Cf&#xA;p2uK0pCtEmCUWK0zKnOSsFwmOfCq5P87VD8AAAD//wMAUEsBAi0AFAAGAAgAAAAhALaDOJL+AAAA&#xA;4QEAABMAAAAAAAAAAAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEA&#xA;OP0h/9YAAACUAQAACwAAAAAAAAAAAAAAAAAvAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEA&#xA;lDVSkIMCAAByBQAADgAAAAAAAAAAAAAAAAAuAgAAZHJzL2Uyb0RvYy54bWxQSwECLQAUAAYACAAA&#xA;ACEAY8912+EAAAALAQAADwAAAAAAAAAAAAAAAADdBAAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAAE&#xA;AAQA8wAAAOsFAAAAAA==&#xA;" fillcolor="#4472c4 [3204]" stroked="f"><v:shadow on="t" color="black" opacity="41287f" offset="0,1.5pt"/><v:textbox><w:txbxContent><w:p w14:paraId="70180F89" w14:textId="3DF3CBB2" w:rsidR="00347169" w:rsidRPr="00347169" w:rsidRDefault="00347169" w:rsidP="00347169"><w:pPr><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="44"/><w:szCs w:val="40"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r w:rsidRPr="00347169"><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="44"/><w:szCs w:val="40"/><w:lang w:val="es-ES"/></w:rPr><w:t>AYUDA</w:t></w:r></w:p></w:txbxContent></v:textbox></v:roundrect></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r w:rsidR="00620EB4"><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="46497C20" wp14:editId="58C38527"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>218961</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>2318233</wp:posOffset></wp:positionV><wp:extent cx="1730707" cy="589221"/><wp:effectExtent l="57150" t="38100" r="60325" b="78105"/><wp:wrapNone/><wp:docPr id="2" name="Cuadro de texto 2"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1730707" cy="589221"/></a:xfrm><a:prstGeom prst="roundRect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FF0000"/></a:solidFill><a:ln/></wps:spPr><wps:style><a:lnRef idx="0"><a:schemeClr val="accent4"/></a:lnRef><a:fillRef idx="3"><a:schemeClr val="accent4"/></a:fillRef><a:effectRef idx="3"><a:schemeClr val="accent4"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:txbx><w:txbxContent><w:p w14:paraId="15CF4A94" w14:textId="0084F2FF" w:rsidR="00620EB4" w:rsidRPr="00347169" w:rsidRDefault="00347169"><w:pPr><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="48"/><w:szCs w:val="44"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="48"/><w:szCs w:val="44"/><w:lang w:val="es-ES"/></w:rPr><w:t>PUNTAJE</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:roundrect w14:anchorId="46497C20" id="Cuadro de texto 2" o:spid="_x0000_s1028" style="position:absolute;margin-left:17.25pt;margin-top:182.55pt;width:136.3pt;height:46.4pt;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:top" arcsize="10923f" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQBuG3RijwIAAHYFAAAOAAAAZHJzL2Uyb0RvYy54bWysVEtvGyEQvlfqf0Dcm11vnDqxso5cR64q&#xA;RU2UpMoZs2AjsQwF7F3313dgH0nT9NCqe2CB+WaY+eZxedXWmhyE8wpMSScnOSXCcKiU2Zb02+P6&#xA;wzklPjBTMQ1GlPQoPL1avH932di5KGAHuhKOoBHj540t6S4EO88yz3eiZv4ErDAolOBqFvDotlnl&#xA;WIPWa50Vef4xa8BV1gEX3uPtdSeki2RfSsHDrZReBKJLir6FtLq0buKaLS7ZfOuY3Sneu8H+wYua&#xA;KYOPjqauWWBk79RvpmrFHXiQ4YRDnYGUiosUA0YzyV9F87BjVqRYkBxvR5r8/zPLvx7uHFFVSQtK&#xA;DKsxRas9qxyQSpAg2gCkiCQ11s8R+2ARHdpP0GKyh3uPlzH2Vro6/jEqgnKk+zhSjJYIj0qz03yW&#xA;zyjhKDs7vyiKZCZ71rbOh88CahI3JXWwN9U95jHRyw43PqA7iB9w8UUPWlVrpXU6uO1mpR05MMz5&#xA;ep3jFz1FlV9g2sTLGFfnf9qFoxbRhjb3QiItKYxkNBakGM0yzoUJ095uQkeURBdGxdPkcarkPyn2&#xA;+KgqUrH+jfKokV4GE0blWhlwb72uw8C27PADA13ckYLQbtpUD2N6N1AdMesOuubxlq8VpuaG+XDH&#xA;HHYLJhonQLjFRWpoSgr9jpIduB9v3Uc8FjFKKWmw+0rqv++ZE5ToLwbL+2IyncZ2TYfp2azAg3sp&#xA;2byUmH29Akz2BGeN5Wkb8UEPW+mgfsJBsYyvoogZjm+XNAzbVehmAg4aLpbLBMIGtSzcmAfLo+nI&#xA;cqy5x/aJOdtXZ+yQrzD0KZu/qs8OGzUNLPcBpErFG3nuWO35x+ZOBdoPojg9Xp4T6nlcLn4CAAD/&#xA;/wMAUEsDBBQABgAIAAAAIQCqxZKt3gAAAAoBAAAPAAAAZHJzL2Rvd25yZXYueG1sTI9Nb8IwDIbv&#xA;k/YfIk/abSQMCqxriqZ90F3XwT00XlvROFUToPz7eSd2si0/ev04W4+uEyccQutJw3SiQCBV3rZU&#xA;a9h+fzysQIRoyJrOE2q4YIB1fnuTmdT6M33hqYy14BAKqdHQxNinUoaqQWfCxPdIvPvxgzORx6GW&#xA;djBnDnedfFRqIZ1piS80psfXBqtDeXQa+k17qcZS7Qr8NO/F4a1YbVSh9f3d+PIMIuIYrzD86bM6&#xA;5Oy090eyQXQaZvOESa6LZAqCgZlacrPXME+WTyDzTP5/If8FAAD//wMAUEsBAi0AFAAGAAgAAAAh&#xA;ALaDOJL+AAAA4QEAABMAAAAAAAAAAAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAU&#xA;AAYACAAAACEAOP0h/9YAAACUAQAACwAAAAAAAAAAAAAAAAAvAQAAX3JlbHMvLnJlbHNQSwECLQAU&#xA;AAYACAAAACEAbht0Yo8CAAB2BQAADgAAAAAAAAAAAAAAAAAuAgAAZHJzL2Uyb0RvYy54bWxQSwEC&#xA;LQAUAAYACAAAACEAqsWSrd4AAAAKAQAADwAAAAAAAAAAAAAAAADpBAAAZHJzL2Rvd25yZXYueG1s&#xA;UEsFBgAAAAAEAAQA8wAAAPQFAAAAAA==&#xA;" fillcolor="red" stroked="f"><v:shadow on="t" color="black" opacity="41287f" offset="0,1.5pt"/><v:textbox><w:txbxContent><w:p w14:paraId="15CF4A94" w14:textId="0084F2FF" w:rsidR="00620EB4" w:rsidRPr="00347169" w:rsidRDefault="00347169"><w:pPr><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="48"/><w:szCs w:val="44"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Comic Sans MS" w:hAnsi="Comic Sans MS"/><w:sz w:val="48"/><w:szCs w:val="44"/><w:lang w:val="es-ES"/></w:rPr><w:t>PUNTAJE</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p></w:txbxContent></v:textbox></v:roundrect></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r w:rsidR="00620EB4"><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="05E77584" wp14:editId="2DA79F9A"><wp:extent cx="2218644" cy="903768"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1" name="Imagen 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 2"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId4"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="2247322" cy="915450"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@
$rng.InsertXML($frag)
Write-Output "done"
